$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 60, shifting existing rows 60-92 down to 61-93.
$ws.Rows.Item(60).EntireRow.Insert()

# Populate the new row 60 with the new weekly price record.
$ws.Range("A60").Value = 10
$ws.Range("B60").Value = "Vega Modelo de Temuco"
$ws.Range("C60").Value = "La Araucanía"
$ws.Range("D60").Value = 45040
$ws.Range("E60").Value = 9
$ws.Range("F60").Value = "Fruta"
$ws.Range("G60").Value = 100107
$ws.Range("H60").Value = "Otros"
$ws.Range("I60").Value = 100107011
$ws.Range("J60").Value = "Tuna"
$ws.Range("K60").Value = "Sin especificar"
$ws.Range("L60").Value = "Especial"
$ws.Range("M60").Value = 65
$ws.Range("N60").Value = 26000
$ws.Range("O60").Value = 26000
$ws.Range("P60").Value = 26000
$ws.Range("Q60").Value = "$/caja 16 kilos"
$ws.Range("R60").Value = "Provincia de Los Andes"
$ws.Range("S60").Value = 1625
$ws.Range("T60").Value = 16
